# NIT-9009788538.xlsx — "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The account-statement sheet used to list 9 workers / 10 period rows (rows 16-25).
# The edit keeps only worker JAIRO DE JESUS CABARCAS ANAYA (rows 16-17, periods
# 2208 and 2209, now shown in swapped order) and drops the other eight workers'
# rows entirely. It also updates the summary totals (Valor Mora, Cant.
# Trabajadores, Cant. Periodos) to reflect the smaller data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Swap the two period/value pairs shown for JAIRO, so period 2209
#        (value 40000) is listed first (row 16) and 2208 (value 37333) second
#        (row 17) -------------------------------------------------------
$periodRow16 = $ws.Range("E16").Value2
$valueRow16  = $ws.Range("F16").Value2
$periodRow17 = $ws.Range("E17").Value2
$valueRow17  = $ws.Range("F17").Value2

$ws.Range("E16").Value2 = $periodRow17
$ws.Range("F16").Value2 = $valueRow17
$ws.Range("E17").Value2 = $periodRow16
$ws.Range("F17").Value2 = $valueRow16

# --- 2. Row 17 becomes the LAST row of the (now two-row) data table, so it
#        should carry the bottom/closing border formatting that used to
#        belong to row 25 (the previous last row). Copy formats only, the
#        values already set above are left untouched. ---------------------
$ws.Range("B25:J25").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3. Remove the rows for the other eight workers (old rows 18-25); this
#        shifts the trailing signature rows (old 30/31) up to 22/23. -------
$ws.Range("18:25").Delete()

# --- 4. Refresh the summary figures for the reduced data set. --------------
$ws.Range("E11").Value2 = 77333   # VALOR MORA total
$ws.Range("C13").Value2 = 1       # Cant. Trabajadores
$ws.Range("F13").Value2 = 2       # Cant. Periodos
